# Apply weekly work report update:
#  - Update "Report Generated On" timestamp in D5
#  - Zero out Total Billed Amount (C8) and all daily "Pricing" (H column) values,
#    including each day's TOTAL row, reflecting a no-charge / corrected billing run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work Report")

# Update the "Report Generated On" timestamp
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"

# Zero out the Total Billed Amount summary figure
$ws.Range("C8").Value = 0

# Zero out all Pricing (H column) values across the four day sections,
# including each day's TOTAL row (H95, H102, H108, H114)
$ws.Range("H16:H95").Value = 0
$ws.Range("H100:H102").Value = 0
$ws.Range("H107:H108").Value = 0
$ws.Range("H113:H114").Value = 0
